$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '257.21'
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = '4.82%'
$ws.Cells.Item(2, 5).Style = "Normal"

$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = '-2.87%'
$ws.Cells.Item(3, 5).Style = "Normal"

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '5.207'
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = '-1.46%'
$ws.Cells.Item(4, 5).Style = "Normal"

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '0.05914'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = '3.54%'
$ws.Cells.Item(5, 5).Style = "Normal"

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '6.668'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = '0.45%'
$ws.Cells.Item(6, 5).Style = "Normal"

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.8659'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = '1.82%'
$ws.Cells.Item(7, 5).Style = "Normal"

$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = '14.65%'
$ws.Cells.Item(8, 5).Style = "Normal"

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.1418'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = '2.66%'
$ws.Cells.Item(9, 5).Style = "Normal"

$ws.Cells.Item(10, 2).Value = 'LiechtensteinCryptoassetsExchange'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.03580'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = '7.41%'
$ws.Cells.Item(10, 5).Style = "Normal"

$ws.Cells.Item(11, 2).Value = 'MandalaExchangeToken'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.07186'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = '1.36%'
$ws.Cells.Item(11, 5).Style = "Normal"

$ws.Cells.Item(12, 2).Value = 'BitrueCoin'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.03156'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = '0.38%'
$ws.Cells.Item(12, 5).Style = "Normal"

$ws.Cells.Item(13, 2).Value = 'BitMartToken'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.09230'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = '0.11%'
$ws.Cells.Item(13, 5).Style = "Normal"

$ws.Cells.Item(14, 2).Value = 'BitForexToken'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.001537'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = '0.93%'
$ws.Cells.Item(14, 5).Style = "Normal"

$ws.Cells.Item(15, 2).Value = 'One'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.0006078'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = '1.69%'
$ws.Cells.Item(15, 5).Style = "Normal"

$ws.Cells.Item(16, 2).Value = 'TigerCash'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.005963'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = '-1.57%'
$ws.Cells.Item(16, 5).Style = "Normal"

$ws.Cells.Item(17, 2).Value = 'LEO'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '3.479'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = '-0.39%'
$ws.Cells.Item(17, 5).Style = "Normal"

$ws.Cells.Item(18, 2).Value = 'GateToken'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '3.267'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = '1.77%'
$ws.Cells.Item(18, 5).Style = "Normal"

$ws.Cells.Item(19, 2).Value = 'BTSEToken'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '2.224'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = '1.63%'
$ws.Cells.Item(19, 5).Style = "Normal"

$ws.Cells.Item(20, 2).Value = 'BitpandaEcosystemToken'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.3151'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = '-0.52%'
$ws.Cells.Item(20, 5).Style = "Normal"

$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = '0.03%'
$ws.Cells.Item(21, 5).Style = "Normal"

$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = '-0.56%'
$ws.Cells.Item(22, 5).Style = "Normal"

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.04173'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = '2.44%'
$ws.Cells.Item(23, 5).Style = "Normal"

$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = '1.47%'
$ws.Cells.Item(24, 5).Style = "Normal"

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.001219'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = '-0.26%'
$ws.Cells.Item(25, 5).Style = "Normal"

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.004515'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = '8.70%'
$ws.Cells.Item(26, 5).Style = "Normal"

$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = '-0.01%'
$ws.Cells.Item(27, 5).Style = "Normal"

$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = '33.78%'
$ws.Cells.Item(28, 5).Style = "Normal"

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.03823'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = '0.87%'
$ws.Cells.Item(40, 5).Style = "Normal"

$ws.Cells.Item(41, 2).Value = 'BKEXToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.1103'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = '3.33%'
$ws.Cells.Item(41, 5).Style = "Normal"

$ws.Cells.Item(42, 2).Value = 'KickToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.003978'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = '-23.54%'
$ws.Cells.Item(42, 5).Style = "Normal"

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.002460'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = '3.79%'
$ws.Cells.Item(43, 5).Style = "Normal"

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.01070'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = '12.90%'
$ws.Cells.Item(44, 5).Style = "Normal"

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.00005433'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = '2.78%'
$ws.Cells.Item(45, 5).Style = "Normal"

$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = '0.02%'
$ws.Cells.Item(46, 5).Style = "Normal"

$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = '22.47%'
$ws.Cells.Item(47, 5).Style = "Normal"

$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = '-4.18%'
$ws.Cells.Item(48, 5).Style = "Normal"

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.00002100'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = '0.02%'
$ws.Cells.Item(49, 5).Style = "Normal"

$ws.Cells.Item(50, 5).NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = '0.02%'
$ws.Cells.Item(50, 5).Style = "Normal"
